{"js": "// The edit relocates several blocks of paragraph text to different\n// paragraphs while leaving paragraph styles / run formatting (bold,\n// italic) anchored in place. We therefore perform the change as a\n// sequence of targeted text replacements, scoped per paragraph so the\n// (sometimes-repeated) search strings cannot match the wrong location.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- \"Objetivos\" section (paragraphs 5 & 6) ------------------------\n// Gets the \"Programa resumido\" summary text (PT + EN) that used to live\n// further down the document.\nconst pObjetivosPt = paragraphs.items[5];\npObjetivosPt.getRange().insertText(\n  \"Introdu\u00e7\u00e3o \u00e0 programa\u00e7\u00e3o em Python; palavras-chave em Python; rotinas e fun\u00e7\u00f5es; classes; numpy e o conceito de slicing e indexing de arrays; revis\u00e3o de m\u00e9todos num\u00e9ricos usando scipy; gera\u00e7\u00e3o de gr\u00e1ficos e anima\u00e7\u00f5es com a biblioteca matplotlib; cria\u00e7\u00e3o de interfaces gr\u00e1ficas com o usu\u00e1rio usando matplotlib.widgets\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst pObjetivosEn = paragraphs.items[6];\npObjetivosEn.getRange().insertText(\n  \"Introduction to Python programming; keywords in Python; routines and functions; classes; numpy and the concept of slicing and indexing arrays; review of numerical methods using scipy; generating graphics and animations with the matplotlib library; creating graphical user interfaces using matplotlib.widgets\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- \"Docente(s) Respons\u00e1vel(eis)\" list (paragraph 8) --------------\n// Its two lines (professor names) are replaced by the old \"Objetivos\"\n// paragraph and the old \"Programa\" paragraph text.\nconst pDocentes = paragraphs.items[8];\n\nlet res = pDocentes.search(\"7290967 - Emerson Gon\u00e7alves de Melo\", { matchCase: true });\nres.load(\"items\");\nawait context.sync();\nres.items[0].insertText(\n  \"Fornecer ao aluno uma introdu\u00e7\u00e3o \u00e0 computa\u00e7\u00e3o cient\u00edfica moderna, usando a linguagem Python e suas bibliotecas num\u00e9ricas e gr\u00e1ficas mais populares: numpy, scipy, matplotlib e pandas. Ao final do curso, o aluno estar\u00e1 capacitado a desenvolver programas complexos, de pequeno e m\u00e9dio porte para solucionar problemas de engenharia que envolvam processamento num\u00e9rico de grandes conjuntos de dados e correlacionar vari\u00e1veis usando m\u00e9todos num\u00e9ricos.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nres = pDocentes.search(\"1176388 - Luiz Tadeu Fernandes Eleno\", { matchCase: true });\nres.load(\"items\");\nawait context.sync();\nres.items[0].insertText(\n  \"\u2022 Introdu\u00e7\u00e3o \u00e0 programa\u00e7\u00e3o em Python \u2022 Instala\u00e7\u00e3o de uma distribui\u00e7\u00e3o Python em Windows e Linux \u2022 Formata\u00e7\u00e3o de arquivos em Python \u2022 Estruturas condicionais \u2022 La\u00e7os de repeti\u00e7\u00e3o de comandos \u2022 Outras palavras-chaves e m\u00e9todos \u2022 Rotinas e fun\u00e7\u00f5es \u2022 C\u00f3digos multifonte e bibliotecas pessoais \u2022 Bibliotecas num\u00e9ricas e gr\u00e1ficas: numpy, scipy e matplotlib \u2022 Programa\u00e7\u00e3o orientada a objeto: classes \u2022 Conceito de objetos e inst\u00e2ncias \u2022 Classes e subclasses\u2022 \u201cArrays\u201d em numpy \u2022 O conceito de array em numpy \u2022 \u201cSlicing\u201d e indexa\u00e7\u00e3o \u2022 Trabalhando com arquivos (entrada e sa\u00edda) \u2022 Gr\u00e1ficos em matplotlib \u2022 A biblioteca matplotlib.pyplot e gr\u00e1ficos em 2D e 3D \u2022 A biblioteca matplotlib.animation para criar gr\u00e1ficos animados. \u2022 Interfaces gr\u00e1ficas com o usu\u00e1rio (Graphical User Interface, GUI) \u2022 Interfaces simples com a biblioteca matplotlib.widgets.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- \"Programa resumido\" section (paragraphs 10 & 11) --------------\n// Becomes the teaching-method text (PT) and the old \"Objetivos\" EN text.\nconst pResumidoPt = paragraphs.items[10];\npResumidoPt.getRange().insertText(\n  \"Aulas expositivas e em laborat\u00f3rio computacional, trabalhos e exerc\u00edcios comentados.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst pResumidoEn = paragraphs.items[11];\npResumidoEn.getRange().insertText(\n  \"Provide the student with an introduction to modern scientific computing, using the Python language and its most popular numerical and graphical libraries: numpy, scipy, matplotlib, and pandas. At the end of the course, the student will be able to develop complex, small and medium-sized programs to solve engineering problems that involve numerical processing of large data sets and correlate variables using numerical methods.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- \"Programa\" section (paragraph 13, PT body) ---------------------\n// Becomes the grading-criteria text. (Paragraph 14, the EN bullet list,\n// is unchanged.)\nconst pProgramaPt = paragraphs.items[13];\npProgramaPt.getRange().insertText(\n  \"M\u00e9dia aritm\u00e9tica de exerc\u00edcios e trabalhos propostos ao longo do curso e uma apresenta\u00e7\u00e3o final de projeto.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---- \"Avalia\u00e7\u00e3o\" list (paragraph 16) --------------------------------\n// The bold \"M\u00e9todo:\" / \"Crit\u00e9rio:\" / \"Norma de recupera\u00e7\u00e3o:\" labels stay\n// put; only the values after them change. Replace the last (unique at\n// the time) value first to avoid ambiguous matches against repeated text.\nconst pAvaliacao = paragraphs.items[16];\n\nres = pAvaliacao.search(\"N\u00e3o haver\u00e1 exame de recupera\u00e7\u00e3o.\", { matchCase: true });\nres.load(\"items\");\nawait context.sync();\nres.items[0].insertText(\"7290967 - Emerson Gon\u00e7alves de Melo\", Word.InsertLocation.replace);\nawait context.sync();\n\nres = pAvaliacao.search(\n  \"M\u00e9dia aritm\u00e9tica de exerc\u00edcios e trabalhos propostos ao longo do curso e uma apresenta\u00e7\u00e3o final de projeto.\",\n  { matchCase: true }\n);\nres.load(\"items\");\nawait context.sync();\nconst bibliografiaLines = [\n  \"Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.\",\n  \"Nilo Ney Coutinho Menezes. Introdu\u00e7\u00e3o \u00e0 Programa\u00e7\u00e3o com Python: Algoritmos e L\u00f3gica de Programa\u00e7\u00e3o Para Iniciantes, 3a ed, 2019.\",\n  \"Ramalho, L. Python Fluente. O\u2019Reilly-Novatec, 2015\",\n  \"Downey, A. B. Pense em Python. O\u2019Reilly-Novatec, 2016.\",\n  \"STEWART, J. M. Python for scientists. Cambridge University Press, 2014.\",\n  \"TELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.\",\n  \"LUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O\u2019Reilly Media, 2006.\",\n  \"MCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015.\",\n].join(\"\\v\");\nres.items[0].insertText(bibliografiaLines, Word.InsertLocation.replace);\nawait context.sync();\n\nres = pAvaliacao.search(\n  \"Aulas expositivas e em laborat\u00f3rio computacional, trabalhos e exerc\u00edcios comentados.\",\n  { matchCase: true }\n);\nres.load(\"items\");\nawait context.sync();\nres.items[0].insertText(\"N\u00e3o haver\u00e1 exame de recupera\u00e7\u00e3o.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---- \"Bibliografia\" body (paragraph 18) -----------------------------\n// The bibliography list moved up into the \"Avalia\u00e7\u00e3o\" paragraph above;\n// this paragraph now just holds the second professor's name.\nconst pBibliografia = paragraphs.items[18];\npBibliografia.getRange().insertText(\n  \"1176388 - Luiz Tadeu Fernandes Eleno\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# The edit relocates several blocks of paragraph text to different\n# paragraphs while leaving paragraph styles / run formatting (bold,\n# italic) anchored in place. We therefore perform the change as a\n# sequence of targeted text replacements, scoped per paragraph (via a\n# freshly-fetched Range each time) so the sometimes-repeated search\n# strings cannot match the wrong location.\n\n$d = $word.ActiveDocument\n\n# ---- \"Objetivos\" section (paragraphs 6 & 7, 1-based) ----------------\n# Gets the \"Programa resumido\" summary text (PT + EN) that used to live\n# further down the document.\n$d.Paragraphs.Item(6).Range.Text = \"Introdu\u00e7\u00e3o \u00e0 programa\u00e7\u00e3o em Python; palavras-chave em Python; rotinas e fun\u00e7\u00f5es; classes; numpy e o conceito de slicing e indexing de arrays; revis\u00e3o de m\u00e9todos num\u00e9ricos usando scipy; gera\u00e7\u00e3o de gr\u00e1ficos e anima\u00e7\u00f5es com a biblioteca matplotlib; cria\u00e7\u00e3o de interfaces gr\u00e1ficas com o usu\u00e1rio usando matplotlib.widgets\"\n\n$d.Paragraphs.Item(7).Range.Text = \"Introduction to Python programming; keywords in Python; routines and functions; classes; numpy and the concept of slicing and indexing arrays; review of numerical methods using scipy; generating graphics and animations with the matplotlib library; creating graphical user interfaces using matplotlib.widgets\"\n\n# ---- \"Docente(s) Respons\u00e1vel(eis)\" list (paragraph 9, 1-based) ------\n# Its two lines (professor names) are replaced by the old \"Objetivos\"\n# paragraph and the old \"Programa\" paragraph text.\n$d.Paragraphs.Item(9).Range.Find.Execute(\n    \"7290967 - Emerson Gon\u00e7alves de Melo\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Fornecer ao aluno uma introdu\u00e7\u00e3o \u00e0 computa\u00e7\u00e3o cient\u00edfica moderna, usando a linguagem Python e suas bibliotecas num\u00e9ricas e gr\u00e1ficas mais populares: numpy, scipy, matplotlib e pandas. Ao final do curso, o aluno estar\u00e1 capacitado a desenvolver programas complexos, de pequeno e m\u00e9dio porte para solucionar problemas de engenharia que envolvam processamento num\u00e9rico de grandes conjuntos de dados e correlacionar vari\u00e1veis usando m\u00e9todos num\u00e9ricos.\",\n    2) | Out-Null\n\n$d.Paragraphs.Item(9).Range.Find.Execute(\n    \"1176388 - Luiz Tadeu Fernandes Eleno\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\u2022 Introdu\u00e7\u00e3o \u00e0 programa\u00e7\u00e3o em Python \u2022 Instala\u00e7\u00e3o de uma distribui\u00e7\u00e3o Python em Windows e Linux \u2022 Formata\u00e7\u00e3o de arquivos em Python \u2022 Estruturas condicionais \u2022 La\u00e7os de repeti\u00e7\u00e3o de comandos \u2022 Outras palavras-chaves e m\u00e9todos \u2022 Rotinas e fun\u00e7\u00f5es \u2022 C\u00f3digos multifonte e bibliotecas pessoais \u2022 Bibliotecas num\u00e9ricas e gr\u00e1ficas: numpy, scipy e matplotlib \u2022 Programa\u00e7\u00e3o orientada a objeto: classes \u2022 Conceito de objetos e inst\u00e2ncias \u2022 Classes e subclasses\u2022 \u201cArrays\u201d em numpy \u2022 O conceito de array em numpy \u2022 \u201cSlicing\u201d e indexa\u00e7\u00e3o \u2022 Trabalhando com arquivos (entrada e sa\u00edda) \u2022 Gr\u00e1ficos em matplotlib \u2022 A biblioteca matplotlib.pyplot e gr\u00e1ficos em 2D e 3D \u2022 A biblioteca matplotlib.animation para criar gr\u00e1ficos animados. \u2022 Interfaces gr\u00e1ficas com o usu\u00e1rio (Graphical User Interface, GUI) \u2022 Interfaces simples com a biblioteca matplotlib.widgets.\",\n    2) | Out-Null\n\n# ---- \"Programa resumido\" section (paragraphs 11 & 12, 1-based) ------\n# Becomes the teaching-method text (PT) and the old \"Objetivos\" EN text.\n$d.Paragraphs.Item(11).Range.Text = \"Aulas expositivas e em laborat\u00f3rio computacional, trabalhos e exerc\u00edcios comentados.\"\n\n$d.Paragraphs.Item(12).Range.Text = \"Provide the student with an introduction to modern scientific computing, using the Python language and its most popular numerical and graphical libraries: numpy, scipy, matplotlib, and pandas. At the end of the course, the student will be able to develop complex, small and medium-sized programs to solve engineering problems that involve numerical processing of large data sets and correlate variables using numerical methods.\"\n\n# ---- \"Programa\" section (paragraph 14, 1-based, PT body) ------------\n# Becomes the grading-criteria text. (Paragraph 15, the EN bullet list,\n# is unchanged.)\n$d.Paragraphs.Item(14).Range.Text = \"M\u00e9dia aritm\u00e9tica de exerc\u00edcios e trabalhos propostos ao longo do curso e uma apresenta\u00e7\u00e3o final de projeto.\"\n\n# ---- \"Avalia\u00e7\u00e3o\" list (paragraph 17, 1-based) ------------------------\n# The bold \"M\u00e9todo:\" / \"Crit\u00e9rio:\" / \"Norma de recupera\u00e7\u00e3o:\" labels stay\n# put; only the values after them change. Replace the last (unique at\n# the time) value first to avoid ambiguous matches against repeated text.\n$d.Paragraphs.Item(17).Range.Find.Execute(\n    \"N\u00e3o haver\u00e1 exame de recupera\u00e7\u00e3o.\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"7290967 - Emerson Gon\u00e7alves de Melo\",\n    2) | Out-Null\n\n$bibliografiaLines = \"Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.^lNilo Ney Coutinho Menezes. Introdu\u00e7\u00e3o \u00e0 Programa\u00e7\u00e3o com Python: Algoritmos e L\u00f3gica de Programa\u00e7\u00e3o Para Iniciantes, 3a ed, 2019.^lRamalho, L. Python Fluente. O\u2019Reilly-Novatec, 2015^lDowney, A. B. Pense em Python. O\u2019Reilly-Novatec, 2016.^lSTEWART, J. M. Python for scientists. Cambridge University Press, 2014.^lTELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.^lLUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O\u2019Reilly Media, 2006.^lMCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015.\"\n\n$d.Paragraphs.Item(17).Range.Find.Execute(\n    \"M\u00e9dia aritm\u00e9tica de exerc\u00edcios e trabalhos propostos ao longo do curso e uma apresenta\u00e7\u00e3o final de projeto.\", $false, $false, $false, $false, $false, $true, 1, $false,\n    $bibliografiaLines,\n    2) | Out-Null\n\n$d.Paragraphs.Item(17).Range.Find.Execute(\n    \"Aulas expositivas e em laborat\u00f3rio computacional, trabalhos e exerc\u00edcios comentados.\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"N\u00e3o haver\u00e1 exame de recupera\u00e7\u00e3o.\",\n    2) | Out-Null\n\n# ---- \"Bibliografia\" body (paragraph 19, 1-based) ---------------------\n# The bibliography list moved up into the \"Avalia\u00e7\u00e3o\" paragraph above;\n# this paragraph now just holds the second professor's name.\n$d.Paragraphs.Item(19).Range.Text = \"1176388 - Luiz Tadeu Fernandes Eleno\"\n"}
